# Updates the "广州-漫展信息" workbook output to the scraper run at 456a3b4.
#
# Sheet 1 "展览" (exhibitions): refresh "想去人数" (want-to-go) counters for
#   a handful of still-open listings (no rows added/removed).
# Sheet 2 "演出" (performances): the two oldest listings (2024-05-25) have
#   expired and are dropped from the feed; later rows shift up and the
#   index column is renumbered.
# Sheet 3 "本地生活" (local life): unchanged.
# Sheet 4 "全部类型" (all types / union of the above): same two expired
#   listings dropped, plus the same counter refresh as sheet 1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 展览 - refresh "想去人数" (column F) for a set of rows.
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    4  = 367
    5  = 1784
    7  = 1462
    8  = 837
    9  = 364
    10 = 712
    11 = 13030
    12 = 12947
    13 = 973
    14 = 753
    16 = 538
    18 = 612
    19 = 2032
    20 = 46
    22 = 23
    24 = 142
    25 = 262
    26 = 710
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value2 = $exhibitUpdates[$row]
}

# Keep a link -> new-count map handy; sheet 4 carries the same rows and
# needs the identical refreshed counters.
$linkToNewCount = @{}
foreach ($row in $exhibitUpdates.Keys) {
    $link = $wsExhibit.Cells.Item($row, 8).Value2
    $linkToNewCount[$link] = $exhibitUpdates[$row]
}

# ---------------------------------------------------------------------
# 2) 演出 - drop the two expired 2024-05-25 rows (rows 2 and 3), shift
#    everything else up, and renumber the index column.
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("A2:A3").EntireRow.Delete() | Out-Null

$lastRow = $wsShow.Cells.Item($wsShow.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $wsShow.Cells.Item($r, 1).Value2 = $r - 1
}

# ---------------------------------------------------------------------
# 3) 本地生活 - no changes.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 4) 全部类型 - same two expired rows removed (found by matching the
#    Link column, since their position differs from sheet 2), plus the
#    refreshed "想去人数" counters that came from 展览.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$expiredLinks = @(
    "https://show.bilibili.com/platform/detail.html?id=85484",
    "https://show.bilibili.com/platform/detail.html?id=83327"
)

$lastRowAll = $wsAll.Cells.Item($wsAll.Rows.Count, 1).End(-4162).Row

# Apply the refreshed counters first (rows keep their original position
# for this pass).
for ($r = 2; $r -le $lastRowAll; $r++) {
    $link = $wsAll.Cells.Item($r, 8).Value2
    if ($linkToNewCount.ContainsKey($link)) {
        $wsAll.Cells.Item($r, 6).Value2 = $linkToNewCount[$link]
    }
}

# Now find and delete the rows for the two expired listings (walk from
# the bottom up so row numbers of not-yet-processed rows stay valid).
for ($r = $lastRowAll; $r -ge 2; $r--) {
    $link = $wsAll.Cells.Item($r, 8).Value2
    if ($expiredLinks -contains $link) {
        $wsAll.Rows.Item($r).Delete() | Out-Null
    }
}

$lastRowAll = $wsAll.Cells.Item($wsAll.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowAll; $r++) {
    $wsAll.Cells.Item($r, 1).Value2 = $r - 1
}
